# Swap the contents of row 13 <-> row 14, and row 23 <-> row 24.
# (Row numbers / positions stay the same; only the record data moves.)
#
# Note: plain Value assignment on this runtime sometimes mis-handles
# reading back a cell's Value property, so Value2 is used for reads.
# Also, columns that store numeric-looking text (e.g. "1", "2") must be
# written with a leading apostrophe so they remain text, matching the
# original inlineStr string cells instead of turning into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 13 <-> 14 (Blasippa / Spillkraka records)
# ---------------------------------------------------------------------

$ws.Range("A13").Value = 112128596
$ws.Range("B13").Value = 56414
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 100049
$ws.Range("F13").Value = "Spillkr" + [char]0x00E5 + "ka"
$ws.Range("G13").Value = "Dryocopus martius"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("I13").Value = "'1"
$ws.Range("J13").ClearContents()
$ws.Range("M13").Value = "f" + [char]0x00E4 + "rska sp" + [char]0x00E5 + "r"
$ws.Range("Q13").Value = 654869.2220899891
$ws.Range("R13").Value = 6675889.891051496
$ws.Range("Z13").Value = "13:35"
$ws.Range("AB13").Value = "13:35"

$ws.Range("A14").Value = 112129086
$ws.Range("B14").Value = 98535
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 222498
$ws.Range("F14").Value = "Bl" + [char]0x00E5 + "sippa"
$ws.Range("G14").Value = "Hepatica nobilis"
$ws.Range("H14").Value = "Schreb."
$ws.Range("I14").ClearContents()
$ws.Range("J14").Value = ""
$ws.Range("M14").ClearContents()
$ws.Range("Q14").Value = 654854.0122225143
$ws.Range("R14").Value = 6675895.72043351
$ws.Range("Z14").Value = "13:11"
$ws.Range("AB14").Value = "13:11"

# ---------------------------------------------------------------------
# Rows 23 <-> 24 (Fjallig taggsvamp s.str. records)
# ---------------------------------------------------------------------

$ws.Range("A23").Value = 112129057
$ws.Range("I23").Value = "'1"
$ws.Range("Q23").Value = 654969.2522851203
$ws.Range("R23").Value = 6675730.300420964
$ws.Range("Z23").Value = "14:02"
$ws.Range("AB23").Value = "14:02"

$ws.Range("A24").Value = 112129050
$ws.Range("I24").Value = "'2"
$ws.Range("Q24").Value = 654923.528877756
$ws.Range("R24").Value = 6675761.7273869
$ws.Range("Z24").Value = "12:20"
$ws.Range("AB24").Value = "12:20"
